$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores its figures as text (the source file uses
# inline strings, e.g. "42.30", "0.320"), so force text format before
# writing the new values to avoid Excel's automatic number conversion,
# which would both change the cell type and drop meaningful trailing
# zeros (e.g. "42.30" -> 42.3, "0.320" -> 0.32).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.537.91'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '3.440.41'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '408.01'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = '133.74'
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.686'
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  -4.62%  '
$ws.Range("D11").Value = '42.29'
$ws.Range("E11").Value = '  -1.75%  '
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").Value = '8.47'
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("D14").Value = '19.96'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = '3.371.45'
$ws.Range("E15").Value = '  -2.81%  '
$ws.Range("D16").Value = '62.557.76'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '11.38'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("D18").Value = '1.03'
$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").Value = '3.19'
$ws.Range("E20").Value = '  -5.46%  '
$ws.Range("D21").Value = '84.18'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").Value = '315.43'
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("D23").Value = '12.94'
$ws.Range("E23").Value = '  -1.80%  '
$ws.Range("D24").Value = '3.16'
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("E25").Value = '  +8.35%  '
$ws.Range("D26").Value = '29.78'
$ws.Range("E26").Value = '  -2.32%  '
$ws.Range("D27").Value = '8.27'
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  +3.04%  '
$ws.Range("D29").Value = '7.57'
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("E30").Value = '  -3.24%  '
$ws.Range("D31").Value = '0.115'
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("D32").Value = '42.30'
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("D35").Value = '0.0486'
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("D36").Value = '51.44'
$ws.Range("E36").Value = '  -2.25%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("E38").Value = '  -5.68%  '
$ws.Range("D39").Value = '2.96'
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D40").Value = '0.320'
$ws.Range("E40").Value = '  +11.38%  '
$ws.Range("D41").Value = '1.99'
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '137.90'
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.125'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = '4.03'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("D45").Value = '16.86'
$ws.Range("E45").Value = '  -4.50%  '
$ws.Range("D46").Value = '2.22'
$ws.Range("E46").Value = '  -1.55%  '
$ws.Range("D47").Value = '21.42'
$ws.Range("E47").Value = '  -5.42%  '
$ws.Range("D48").Value = '2.129.68'
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").Value = '2.33'
$ws.Range("E49").Value = '  -3.70%  '
$ws.Range("B50").Value = 'Fetch.AI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D50").Value = '1.75'
$ws.Range("E50").Value = '  +21.89%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '1.92'
$ws.Range("E51").Value = '  +2.17%  '
